# "broke up 'programing' into programig for Server and Programming of the Arduino"
#
# The original shared "Programmierung" / "4x36h" pair (row 5, column J/L) is
# replaced by two distinct tasks:
#   - J5/L5 -> "Programmierung Arduino & Sensoren" / "4x18h"
#   - J7/L7 -> "Programmierung Serverseitig" / "4x18h"
# This pushes the rest of the J/K/L "Testbetrieb" task chain down by one slot
# (J7->J9->J11->J13->J15), and a brand new J15/K15/L15 entry is created for
# "Sensorkonfigurationen erstellen".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Fix up formatting for the row-15 block (J15:L15) before writing values,
#     by cloning the formats from the row-13 block immediately above it,
#     which already carries the correct visual style for this task chain. ---
$ws.Range("J13").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("K13").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("L13").Copy()
$ws.Range("L15").PasteSpecial(-4122)

# --- New blank placeholder cells on row 14 (J14/L14), matching the style of
#     the existing blank placeholder H14. ---
$ws.Range("H14").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("L14").PasteSpecial(-4122)

# --- Cells that become entirely blank (no value, no explicit style). ---
$ws.Range("K8").Clear()
$ws.Range("L8").Clear()
$ws.Range("J10").Clear()
$ws.Range("K12").Clear()
$ws.Range("K14").Clear()

# --- Programmierung split: row 5 becomes the Arduino/Sensor task, row 7
#     becomes the server-side task, each with its own 4x18h effort. (Effort
#     cells are written first so "4x18h" is registered as a shared string
#     ahead of the task-name strings.) ---
$ws.Range("L5").Value = "4x18h"
$ws.Range("L7").Value = "4x18h"
$ws.Range("J5").Value = "Programmierung Arduino & Sensoren"
$ws.Range("J7").Value = "Programmierung Serverseitig"

# --- Remaining Testbetrieb tasks each shift down one slot in the chain. ---
$ws.Range("J9").Value = "Testszenarien kreieren, eventuelle Betriebserlaubnisse beantragen"
$ws.Range("L9").Value = "2h"

$ws.Range("E11").Value = 23
$ws.Range("J11").Value = "Korrekturschleifen planen und durchführen"
$ws.Range("L11").Value = "12h"

$ws.Range("J13").Value = "Dokumentation erstellen"
$ws.Range("L13").Value = "8h"

$ws.Range("J15").Value = "Sensorkonfigurationen erstellen"
$ws.Range("K15").Value = 45
$ws.Range("L15").Value = "6h"

# --- Selection moved to M16 in the saved view. ---
$ws.Activate()
$ws.Range("M16").Select()
